# Update "想去人数" (want-to-go count) figures in the F column of the
# "展览" and "全部类型" worksheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 120
$ws1.Range("F3").Value  = 312
$ws1.Range("F4").Value  = 63
$ws1.Range("F5").Value  = 804
$ws1.Range("F7").Value  = 2125
$ws1.Range("F8").Value  = 277
$ws1.Range("F9").Value  = 105
$ws1.Range("F10").Value = 4704
$ws1.Range("F11").Value = 6
$ws1.Range("F13").Value = 294
$ws1.Range("F15").Value = 22
$ws1.Range("F16").Value = 159
$ws1.Range("F17").Value = 32
$ws1.Range("F19").Value = 105
$ws1.Range("F20").Value = 3618
$ws1.Range("F21").Value = 203
$ws1.Range("F22").Value = 586
$ws1.Range("F25").Value = 94
$ws1.Range("F26").Value = 107
$ws1.Range("F29").Value = 77
$ws1.Range("F30").Value = 216
$ws1.Range("F32").Value = 787
$ws1.Range("F33").Value = 2239
$ws1.Range("F34").Value = 411

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 120
$ws4.Range("F3").Value  = 312
$ws4.Range("F4").Value  = 63
$ws4.Range("F5").Value  = 805
$ws4.Range("F7").Value  = 2125
$ws4.Range("F8").Value  = 277
$ws4.Range("F9").Value  = 105
$ws4.Range("F10").Value = 4704
$ws4.Range("F11").Value = 6
$ws4.Range("F13").Value = 294
$ws4.Range("F15").Value = 22
$ws4.Range("F16").Value = 159
$ws4.Range("F17").Value = 32
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 3618
$ws4.Range("F21").Value = 203
$ws4.Range("F22").Value = 586
$ws4.Range("F25").Value = 94
$ws4.Range("F26").Value = 107
$ws4.Range("F29").Value = 77
$ws4.Range("F30").Value = 216
$ws4.Range("F33").Value = 788
$ws4.Range("F34").Value = 2239
$ws4.Range("F35").Value = 411
